# Lingualeo alpha: rewrite the word/phrase lines and tidy up the paragraph list.
$d = $word.ActiveDocument

# --- Paragraph 1: "Test1 — [test1] — тест1" -> "t" + "est1 — тест1" ---
$p1 = $d.Paragraphs(1).Range
$p1.Text = "test1 — тест1"
$start = $d.Paragraphs(1).Range.Start
$r = $d.Range($start, $start + 1)
$r.Font.Bold = 1
$r.Font.Bold = 0

# --- Paragraph 2: "Test2 — [test2]" -> "h" + "ello — " + "п" + "ривет" ---
$p2 = $d.Paragraphs(2).Range
$p2.Text = "hello — привет"
$start = $d.Paragraphs(2).Range.Start

$r = $d.Range($start, $start + 1)
$r.Font.Bold = 1
$r.Font.Bold = 0

$r = $d.Range($start + 1, $start + 8)
$r.Font.Bold = 1
$r.Font.Bold = 0

$r = $d.Range($start + 8, $start + 9)
$r.Font.Bold = 1
$r.Font.Bold = 0

# --- Paragraph 3: "Test3 — тест3" -> "guitar — гитара" ---
$p3 = $d.Paragraphs(3).Range
$p3.Text = "guitar — гитара"

# --- Paragraph 4 ("Test4") is dropped entirely ---
$d.Paragraphs(4).Range.Delete()
